# Update "loading_percent" results for the 380 kV case (Case_3_93).
# Columns B, C, D, F, G, H, N for rows 2-25 receive new computed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 13.4252224796261
$ws.Range("C2").Value = 9.38299382964138
$ws.Range("D2").Value = 3.740008472580512
$ws.Range("F2").Value = 21.38989080404852
$ws.Range("G2").Value = 24.44737171983198
$ws.Range("H2").Value = 12.70049582734677
$ws.Range("N2").Value = 15.88000812443902
$ws.Range("B3").Value = 12.77419563848321
$ws.Range("C3").Value = 8.878361777503271
$ws.Range("D3").Value = 3.718775219632284
$ws.Range("F3").Value = 21.28832621055277
$ws.Range("G3").Value = 24.1813638881743
$ws.Range("H3").Value = 12.73580959058132
$ws.Range("N3").Value = 15.93321286040914
$ws.Range("B4").Value = 12.35888583854069
$ws.Range("C4").Value = 8.551748914040825
$ws.Range("D4").Value = 3.705699616513914
$ws.Range("F4").Value = 21.23491910154956
$ws.Range("G4").Value = 24.0304361432583
$ws.Range("H4").Value = 12.76104410354162
$ws.Range("N4").Value = 15.96786108781562
$ws.Range("B5").Value = 12.18594976194884
$ws.Range("C5").Value = 8.414501547623795
$ws.Range("D5").Value = 3.700364075237709
$ws.Range("F5").Value = 21.21542303914675
$ws.Range("G5").Value = 23.97213113720726
$ws.Range("H5").Value = 12.77221658185119
$ws.Range("N5").Value = 15.98247936321711
$ws.Range("B6").Value = 12.15701796004825
$ws.Range("C6").Value = 8.391463300995355
$ws.Range("D6").Value = 3.69947775503784
$ws.Range("F6").Value = 21.21232308832882
$ws.Range("G6").Value = 23.96264502736274
$ws.Range("H6").Value = 12.77412535119219
$ws.Range("N6").Value = 15.98493687174293
$ws.Range("B7").Value = 12.35656819850367
$ws.Range("C7").Value = 8.549914653067617
$ws.Range("D7").Value = 3.705627685113847
$ws.Range("F7").Value = 21.23464697061064
$ws.Range("G7").Value = 24.02963676858449
$ws.Range("H7").Value = 12.76119118446527
$ws.Range("N7").Value = 15.96805621403706
$ws.Range("B8").Value = 13.20410456575486
$ws.Range("C8").Value = 9.212512008141118
$ws.Range("D8").Value = 3.732696701174478
$ws.Range("F8").Value = 21.35302400759488
$ws.Range("G8").Value = 24.35312952386182
$ws.Range("H8").Value = 12.71193257300615
$ws.Range("N8").Value = 15.89794262882672
$ws.Range("B9").Value = 14.73438812377864
$ws.Range("C9").Value = 10.37658822698788
$ws.Range("D9").Value = 3.78536536425418
$ws.Range("F9").Value = 21.65528248750466
$ws.Range("G9").Value = 25.08189063131426
$ws.Range("H9").Value = 12.64368297292283
$ws.Range("N9").Value = 15.77612381791767
$ws.Range("B10").Value = 15.76948291294286
$ws.Range("C10").Value = 11.14713297558025
$ws.Range("D10").Value = 3.823662593337879
$ws.Range("F10").Value = 21.91856934449638
$ws.Range("G10").Value = 25.66891248120654
$ws.Range("H10").Value = 12.61103226318869
$ws.Range("N10").Value = 15.69612528397333
$ws.Range("B11").Value = 16.21961013938946
$ws.Range("C11").Value = 11.47898190254013
$ws.Range("D11").Value = 3.840967347453681
$ws.Range("F11").Value = 22.04688715757819
$ws.Range("G11").Value = 25.94577796014855
$ws.Range("H11").Value = 12.60001683653899
$ws.Range("N11").Value = 15.66178427081102
$ws.Range("B12").Value = 16.38698854826345
$ws.Range("C12").Value = 11.60194189244629
$ws.Range("D12").Value = 3.847500758481068
$ws.Range("F12").Value = 22.09666735895396
$ws.Range("G12").Value = 26.05191093372638
$ws.Range("H12").Value = 12.59640021963897
$ws.Range("N12").Value = 15.64907435995959
$ws.Range("B13").Value = 16.35107867036049
$ws.Range("C13").Value = 11.57558075827703
$ws.Range("D13").Value = 3.846094583014417
$ws.Range("F13").Value = 22.0858940901093
$ws.Range("G13").Value = 26.02899790481458
$ws.Range("H13").Value = 12.59715441085578
$ws.Range("N13").Value = 15.65179858754702
$ws.Range("B14").Value = 16.23344255966404
$ws.Range("C14").Value = 11.48915217420526
$ws.Range("D14").Value = 3.841505261933034
$ws.Range("F14").Value = 22.05095901523309
$ws.Range("G14").Value = 25.95448444166968
$ws.Range("H14").Value = 12.59970816331849
$ws.Range("N14").Value = 15.6607327242338
$ws.Range("B15").Value = 16.16098413229188
$ws.Range("C15").Value = 11.43585957263271
$ws.Range("D15").Value = 3.838691547160188
$ws.Range("F15").Value = 22.02971387794115
$ws.Range("G15").Value = 25.90900710634255
$ws.Range("H15").Value = 12.60134472631871
$ws.Range("N15").Value = 15.66624344814544
$ws.Range("B16").Value = 15.73964037815549
$ws.Range("C16").Value = 11.12506836479734
$ws.Range("D16").Value = 3.822529102762755
$ws.Range("F16").Value = 21.91035211483518
$ws.Range("G16").Value = 25.65100594199647
$ws.Range("H16").Value = 12.61182961921649
$ws.Range("N16").Value = 15.69841076320405
$ws.Range("B17").Value = 15.47577591946293
$ws.Range("C17").Value = 10.92960986061973
$ws.Range("D17").Value = 3.812582099010606
$ws.Range("F17").Value = 21.83928922144806
$ws.Range("G17").Value = 25.49516062907685
$ws.Range("H17").Value = 12.61924677466331
$ws.Range("N17").Value = 15.71866916708336
$ws.Range("B18").Value = 15.32206002604
$ws.Range("C18").Value = 10.81543017477071
$ws.Range("D18").Value = 3.806850009389489
$ws.Range("F18").Value = 21.79922284796908
$ws.Range("G18").Value = 25.40645498353555
$ws.Range("H18").Value = 12.62387398253777
$ws.Range("N18").Value = 15.73051432291391
$ws.Range("B19").Value = 15.26968297615286
$ws.Range("C19").Value = 10.77646975792064
$ws.Range("D19").Value = 3.804907444355812
$ws.Range("F19").Value = 21.7857968550802
$ws.Range("G19").Value = 25.3765846338085
$ws.Range("H19").Value = 12.62550259551732
$ws.Range("N19").Value = 15.73455806599744
$ws.Range("B20").Value = 15.50406707681612
$ws.Range("C20").Value = 10.95059875820282
$ws.Range("D20").Value = 3.813642115785223
$ws.Range("F20").Value = 21.84677073369107
$ws.Range("G20").Value = 25.51165500648562
$ws.Range("H20").Value = 12.61841981617371
$ws.Range("N20").Value = 15.71649264912238
$ws.Range("B21").Value = 16.2680792456169
$ws.Range("C21").Value = 11.51461185131255
$ws.Range("D21").Value = 3.842853808265721
$ws.Range("F21").Value = 22.06118835822755
$ws.Range("G21").Value = 25.97633683169136
$ws.Range("H21").Value = 12.59894298909567
$ws.Range("N21").Value = 15.65810056956697
$ws.Range("B22").Value = 16.74945016706784
$ws.Range("C22").Value = 11.86746098305125
$ws.Range("D22").Value = 3.861830287449507
$ws.Range("F22").Value = 22.20823256023286
$ws.Range("G22").Value = 26.28749426851204
$ws.Range("H22").Value = 12.58944799435576
$ws.Range("N22").Value = 15.62165303266458
$ws.Range("B23").Value = 16.49420269952031
$ws.Range("C23").Value = 11.68058614239797
$ws.Range("D23").Value = 3.851713605525726
$ws.Range("F23").Value = 22.12913406423433
$ws.Range("G23").Value = 26.12078211925893
$ws.Range("H23").Value = 12.59421885699313
$ws.Range("N23").Value = 15.64094903268493
$ws.Range("B24").Value = 15.49128291898255
$ws.Range("C24").Value = 10.94111530775031
$ws.Range("D24").Value = 3.813162923788398
$ws.Range("F24").Value = 21.84338588241669
$ws.Range("G24").Value = 25.50419511088615
$ws.Range("H24").Value = 12.61879255359475
$ws.Range("N24").Value = 15.71747603502675
$ws.Range("B25").Value = 14.33552692493228
$ws.Range("C25").Value = 10.07642338405769
$ws.Range("D25").Value = 3.771176451941348
$ws.Range("F25").Value = 21.56614866440894
$ws.Range("G25").Value = 24.87522337703711
$ws.Range("H25").Value = 12.6590887605638
$ws.Range("N25").Value = 15.80740668800622
